# Auto-generated: apply scheduled market-price refresh to Profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2320.6155
$ws.Range("I62").Value = 2425.3635
$ws.Range("K62").Value = 2425.3635
$ws.Range("M62").Value = -1801.3635
$ws.Range("H65").Value = 2320.6155
$ws.Range("I65").Value = 2425.3635
$ws.Range("K65").Value = 12126.8175
$ws.Range("M65").Value = -9006.817499999999
$ws.Range("H96").Value = 1436.4
$ws.Range("I96").Value = 1231.6666
$ws.Range("J96").Value = 1743.5
$ws.Range("K96").Value = 3694.9998
$ws.Range("L96").Value = 5230.5
$ws.Range("M96").Value = -2321.9998
$ws.Range("N96").Value = -7976.5
$ws.Range("H137").Value = 3335109.2
$ws.Range("I137").Value = 956.1667
$ws.Range("K137").Value = 2868.5001
$ws.Range("M137").Value = -318.5001000000002
$ws.Range("H138").Value = 4979.8135
$ws.Range("I138").Value = 9300.723
$ws.Range("J138").Value = 3082.8293
$ws.Range("K138").Value = 27902.169
$ws.Range("L138").Value = 9248.4879
$ws.Range("M138").Value = -22762.169
$ws.Range("N138").Value = -19528.4879

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 153668.84
$ws.Range("I32").Value = 156681.6
$ws.Range("J32").Value = 55754.5
$ws.Range("K32").Value = 156681.6
$ws.Range("L32").Value = 55754.5
$ws.Range("M32").Value = -156394.6
$ws.Range("N32").Value = -56328.5
$ws.Range("H45").Value = 146956.42
$ws.Range("I45").Value = 146956.42
$ws.Range("K45").Value = 146956.42
$ws.Range("M45").Value = -146579.42
$ws.Range("H122").Value = 2406.6316
$ws.Range("I122").Value = 2106.9375
$ws.Range("J122").Value = 4005
$ws.Range("K122").Value = 6320.8125
$ws.Range("L122").Value = 12015
$ws.Range("M122").Value = -3870.8125
$ws.Range("N122").Value = -16915

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2775.0312
$ws.Range("I31").Value = 1577.5883
$ws.Range("K31").Value = 1577.5883
$ws.Range("M31").Value = -1282.5883
$ws.Range("H34").Value = 2775.0312
$ws.Range("I34").Value = 1577.5883
$ws.Range("K34").Value = 1577.5883
$ws.Range("M34").Value = -1375.5883
$ws.Range("H58").Value = 2613.158
$ws.Range("J58").Value = 3723.7778
$ws.Range("L58").Value = 3723.7778
$ws.Range("N58").Value = -4129.7778
$ws.Range("H94").Value = 1530.625
$ws.Range("I94").Value = 1156
$ws.Range("J94").Value = 1755.4
$ws.Range("K94").Value = 1156
$ws.Range("L94").Value = 1755.4
$ws.Range("M94").Value = -705
$ws.Range("N94").Value = -2657.4
$ws.Range("H105").Value = 2322.8333
$ws.Range("I105").Value = 1412
$ws.Range("K105").Value = 1412
$ws.Range("M105").Value = 335
$ws.Range("H122").Value = 4495.8096
$ws.Range("I122").Value = 2689.611
$ws.Range("K122").Value = 8068.833
$ws.Range("M122").Value = -5618.833
$ws.Range("H136").Value = 2613.158
$ws.Range("J136").Value = 3723.7778
$ws.Range("L136").Value = 11171.3334
$ws.Range("N136").Value = -16271.3334
$ws.Range("H137").Value = 99743.664
$ws.Range("J137").Value = 99743.664
$ws.Range("L137").Value = 99743.664
$ws.Range("N137").Value = -109943.664
$ws.Range("H141").Value = 737599.8
$ws.Range("J141").Value = 737599.8
$ws.Range("L141").Value = 737599.8
$ws.Range("N141").Value = -747959.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 499
$ws.Range("I5").Value = 499
$ws.Range("K5").Value = 1497
$ws.Range("M5").Value = -1385
$ws.Range("H39").Value = 5902.857
$ws.Range("J39").Value = 7149.1816
$ws.Range("L39").Value = 21447.5448
$ws.Range("N39").Value = -22035.5448
$ws.Range("H50").Value = 2008.6666
$ws.Range("I50").Value = 480
$ws.Range("K50").Value = 1440
$ws.Range("M50").Value = -959
$ws.Range("H53").Value = 2008.6666
$ws.Range("I53").Value = 480
$ws.Range("K53").Value = 1440
$ws.Range("M53").Value = -959
$ws.Range("H135").Value = 499
$ws.Range("I135").Value = 499
$ws.Range("K135").Value = 4491
$ws.Range("M135").Value = -1956

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H28").Value = 14998
$ws.Range("J28").Value = 14998
$ws.Range("L28").Value = 14998
$ws.Range("N28").Value = -15382
$ws.Range("H57").Value = 33333
$ws.Range("J57").Value = 33333
$ws.Range("L57").Value = 33333
$ws.Range("N57").Value = -34973
$ws.Range("H98").Value = 22000
$ws.Range("J98").Value = 22000
$ws.Range("L98").Value = 22000
$ws.Range("N98").Value = -27990
$ws.Range("H102").Value = 41668500
$ws.Range("I102").Value = 41668500
$ws.Range("K102").Value = 41668500
$ws.Range("M102").Value = -41666878
$ws.Range("H113").Value = 2783.1765
$ws.Range("I113").Value = 2860
$ws.Range("J113").Value = 2642.3333
$ws.Range("K113").Value = 2860
$ws.Range("L113").Value = 2642.3333
$ws.Range("M113").Value = -690
$ws.Range("N113").Value = -6982.3333
$ws.Range("H122").Value = 7642.2
$ws.Range("I122").Value = 5945.5
$ws.Range("K122").Value = 17836.5
$ws.Range("M122").Value = -15386.5
$ws.Range("H136").Value = 21289.924
$ws.Range("J136").Value = 21289.924
$ws.Range("L136").Value = 63869.772
$ws.Range("N136").Value = -68969.772

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 74984.5
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H22").Value = 2391.9443
$ws.Range("J22").Value = 3162.6667
$ws.Range("L22").Value = 3162.6667
$ws.Range("N22").Value = -3752.6667
$ws.Range("H27").Value = 2391.9443
$ws.Range("J27").Value = 3162.6667
$ws.Range("L27").Value = 3162.6667
$ws.Range("N27").Value = -3376.6667
$ws.Range("H30").Value = 4149.25
$ws.Range("I30").Value = 4149.25
$ws.Range("K30").Value = 4149.25
$ws.Range("M30").Value = -4041.25
$ws.Range("H40").Value = 4335.4375
$ws.Range("I40").Value = 3371.5
$ws.Range("K40").Value = 3371.5
$ws.Range("M40").Value = -3235.5
$ws.Range("H64").Value = 94999.75
$ws.Range("J64").Value = 94999.75
$ws.Range("L64").Value = 94999.75
$ws.Range("N64").Value = -95449.75
$ws.Range("H67").Value = 94999.75
$ws.Range("J67").Value = 94999.75
$ws.Range("L67").Value = 94999.75
$ws.Range("N67").Value = -96559.75
$ws.Range("H122").Value = 4306.256
$ws.Range("I122").Value = 2778.9583
$ws.Range("K122").Value = 8336.874899999999
$ws.Range("M122").Value = -5886.874899999999
$ws.Range("H140").Value = 88099.234
$ws.Range("J140").Value = 88099.234
$ws.Range("L140").Value = 88099.234
$ws.Range("N140").Value = -98459.234

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 25999.5
$ws.Range("J82").Value = 24499
$ws.Range("L82").Value = 24499
$ws.Range("N82").Value = -25265
$ws.Range("H85").Value = 25999.5
$ws.Range("J85").Value = 24499
$ws.Range("L85").Value = 24499
$ws.Range("N85").Value = -27151
$ws.Range("H113").Value = 304.8
$ws.Range("I113").Value = 262.8125
$ws.Range("J113").Value = 472.75
$ws.Range("K113").Value = 788.4375
$ws.Range("L113").Value = 1418.25
$ws.Range("M113").Value = 1381.5625
$ws.Range("N113").Value = -5758.25
$ws.Range("H122").Value = 2257.3809
$ws.Range("I122").Value = 2270.25
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 6810.75
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -4360.75
$ws.Range("N122").Value = -10900
$ws.Range("H140").Value = 87497.5
$ws.Range("J140").Value = 87497.5
$ws.Range("L140").Value = 87497.5
$ws.Range("N140").Value = -97857.5
$ws.Range("H141").Value = 162137
$ws.Range("J141").Value = 162137
$ws.Range("L141").Value = 162137
$ws.Range("N141").Value = -172497
